{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of paragraph index -> [expectedOriginalText, newText]\nconst edits = {\n  2: [\"The article is about how ChatGPT and AI are expected to impact the Saudi workforce productivity by automating mundane tasks, providing tailor-made training programs, and reshaping various sectors like healthcare, transportation, energy, finance, and retail.\", \"The article is about the potential benefits of ChatGPT and AI on Saudi Arabia s workforce productivity, alleviating fears about job losses due to automation by highlighting their role in upskilling employees, boosting productivity, and creating new career opportunities.\"],\n  4: [\"1. ChatGPT and AI can boost productivity in Saudi Arabia by providing tailor-made training programs, access to customized online courses, and fostering collaboration among team members.\\u000b2. Embedding AI into operations requires a holistic approach that defines strategic objectives, advantages, and disadvantages. Understanding operational bottlenecks and how AI can address them is crucial.\\u000b3. The automation of tasks by ChatGPT may lead to job displacement for employees entrusted with mundane and repetitive tasks, necessitating retraining or upskilling for workers.\\u000b4. A positive employee and client experience is vital when implementing AI. Leadership, a futuristic view, and agility in making timely changes are essential for successful integration of AI into organizations.\", \"1. ChatGPT and AI can positively impact Saudi workforce productivity by providing tailor-made training programs, access to customized online courses, fostering collaboration and communication among team members.\\u000b2. Embracing AI will require a strategic approach that defines objectives, advantages, and disadvantages, with an emphasis on addressing operational bottlenecks and choosing the optimal AI tools for the organization's needs.\\u000b3. The automation of tasks by ChatGPT might replace workers entrusted with mundane and repetitive functions, necessitating employees to acquire new skills through retraining or upskilling.\\u000b4. Implementing AI in Saudi companies will require a holistic approach, fostering an optimistic environment for learning and improvement, solid leadership, and agility in making timely changes as needed.\"],\n  6: [\"In the given text, the media frames the public discussion about ChatGPT and AI in general as a tool for boosting productivity and economic development. The metaphor of a wave of change is used to describe the impact of these technologies on the global workforce. Another recurring metaphor is that of an opportunity for organizations and employees to innovate, learn, collaborate, and improve their skills to adapt to an AI-enabled digital world. The text also refers to AI as a tool that can help automate mundane and repetitive tasks to free up time for strategic activities, using the metaphor of a productive and efficient workforce. However, there is also a concern about job displacement due to automation, which is depicted as a potential negative side effect. Overall, the text presents AI, particularly ChatGPT, as a transformative technology that can reshape and emancipate the Saudi workforce, but requires careful management and strategic planning to ensure positive outcomes.\", \"In the text, the public discussion about ChatGPT and AI is framed as a tool for boosting productivity in Saudi Arabia's workforce. The metaphors used to describe this process include:\\u000b\\u000b1. Wave of change (implying the transformative impact of AI technologies on the global workforce)\\u000b2. Opportunity to boost productivity (indicating the potential positive outcomes of embracing innovation and AI)\\u000b3. Embracing innovation (emphasizing the need for organizations to adapt and utilize AI platforms)\\u000b4. Nurturing the right talent with a strong AI culture (implying the importance of fostering a supportive environment for AI implementation)\\u000b5. Reshaping and even emancipating workforce (suggesting that AI can lead to significant changes in the workforce, potentially empowering it)\\u000b6. Tailor-made training programs (indicating personalized learning opportunities provided by AI platforms like ChatGPT)\\u000b7. Automation of tasks (implying the removal of mundane and repetitive tasks from employees)\\u000b8. Acquiring new skills through retraining or upskilling (suggesting that workers will need to adapt to a more AI-enabled digital world)\\u000b9. AI-enabled digital world (emphasizing the growing reliance on artificial intelligence in various sectors)\\u000b10. Strategic roadmap (implying a clear plan for integrating AI into organizations' operations)\\u000b11. Agility (emphasizing the need for adaptability and quick changes as technological advancements occur)\\u000b\\u000bThese metaphors suggest that the media presents ChatGPT and AI as a transformative force, requiring a strategic approach to integration, and necessitating employee adaptation and skill acquisition in order to thrive in an increasingly AI-dependent world.\"],\n  8: [\"Perspectives and aspects that are being widely covered in this text include  1. The positive impact of AI technologies such as ChatGPT on Saudi Arabian workforce productivity. 2. Fear among workers about losing their jobs due to AI, and efforts to alleviate these fears by nurturing the right talent with a strong AI culture. 3. The role of AI in various sectors such as healthcare, transportation, energy, finance, and retail. 4. The benefits of AI for recruitment, hiring, training, development, upskilling, reskilling, talent collaboration, and knowledge management. 5. The potential of AI to increase productivity by automating mundane tasks and freeing up employees to focus on strategic activities. 6. The need to create a positive employee experience through understanding the strategic objectives, advantages, and disadvantages of AI implementation. 7. The importance of leadership, agility, and an optimistic environment for learning and improvement in embedding AI into organizations. 8. The potential impact of AI on Saudi Arabia s public sector operations and service delivery. 9. The concerns about the stability of worker employment due to automation. 10. The need for employees to acquire new skills through retraining or upskilling to remain marketable in an increasingly digital world. Aspects that are being ignored in this text include  1. The specific challenges faced by Saudi Arabia in implementing AI technologies, such as infrastructure limitations, regulatory constraints, or cultural resistance. 2. The ethical considerations related to the use of AI, including issues of data privacy, bias, and accountability. 3. The potential negative consequences of widespread AI adoption on employment, income distribution, and social inequality. 4. The environmental impact of AI technologies, such as energy consumption, e-waste generation, or resource depletion. 5. The potential risks associated with over-reliance on AI, such as algorithmic failures, cybersecurity threats, or dependence on a single vendor. 6. The need for continuous monitoring and evaluation of AI systems to ensure their performance remains optimal and aligned with strategic goals. 7. The need for international cooperation in the development and regulation of AI technologies to prevent a digital divide between developed and developing countries.\", \"The text mainly covers the positive impact of ChatGPT and AI on Saudi workforce productivity, employee training, recruitment, and collaboration, as well as their potential applications in various sectors such as healthcare, transportation, energy, finance, and retail. It also discusses the need to nurture a strong AI culture within organizations, emphasizing innovation, experimentation, learning, and collaboration. The text mentions the concerns about job displacement due to AI but focuses more on how employees can acquire new skills through retraining or upskilling to adapt to an AI-enabled work environment.\\u000b\\u000bAspects that seem to be ignored in this text include the potential negative social impacts of widespread AI adoption, such as increased income inequality and loss of privacy, the ethical considerations around AI decision-making, and the need for regulations and guidelines to ensure fair use of AI technologies. Additionally, the text does not discuss the possible challenges or resistance from employees or organizations in adopting AI, nor does it delve into the long-term implications of AI on the Saudi economy beyond productivity improvements.\"],\n  10: [\"Not mentioned\", \"The Arabic world is not explicitly mentioned in the text, but it can be inferred that the Kingdom (presumably Saudi Arabia) is leveraging AI to boost economic development and increase workforce productivity. The AI platform ChatGPT is highlighted as a useful tool for this purpose. The article discusses how AI, specifically ChatGPT, can impact recruitment, hiring, training, upskilling, reskilling, talent collaboration, knowledge management, and various sectors like healthcare, transportation, energy, finance, and retail. It also mentions the importance of nurturing a strong AI culture to maximize the benefits of AI technologies.\"],\n  12: [\"The article suggests that AI, specifically ChatGPT, can positively impact Saudi workforce productivity by providing customized training programs, collaboration tools, and automating mundane tasks. However, it also warns about potential job displacement due to automation. Embedding AI requires a holistic approach, clear strategic objectives, understanding of operational bottlenecks, and fostering an optimistic learning environment. Leadership, agility, and a futuristic view are crucial for successful implementation.\", \"The article suggests that ChatGPT and AI can boost Saudi Arabia's workforce productivity by automating mundane tasks and providing tailored training programs. However, it emphasizes the need to foster an optimistic learning environment and nurture a strong AI culture to address employees' concerns about job displacement. Embracing innovation can lead to increased efficiency, better positions for employees, and growth opportunities in various sectors like healthcare, transportation, finance, and retail.\"],\n  16: [\"Embedding AI, Arthur D. Little, Saudi Arabias, Kingdoms, Kaspersky, Arab News, Raymond Khoury, Khoury, AI, Khourys, Kasperskys, Arthur D. Little\", \"Kaspersky, AI, Raymond Khoury, Kasperskys, Embedding AI, Kingdoms, Arthur D. Little, Arab News, Khourys, Arthur D. Little, Saudi Arabias, Khoury\"],\n};\n\nfor (const idxStr of Object.keys(edits)) {\n  const idx = Number(idxStr);\n  const [expected, next] = edits[idxStr];\n  const para = paragraphs.items[idx];\n  if (para.text !== expected) {\n    throw new Error(`Paragraph ${idx} text mismatch. Got: ${para.text}`);\n  }\n  para.insertText(next, \"Replace\");\n}\n\nawait context.sync();", "ps1": "# Applies the documented edits to the active Word document using the COM object model.\n# Word represents a manual line break (<w:br/>) as a vertical-tab character (chr 11)\n# in Range.Text, and a paragraph end as a carriage return (chr 13).\n$d = $word.ActiveDocument\n$br = [char]11\n\nfunction Set-ParaText($Document, $Index, $Expected, $NewText) {\n    $para = $Document.Paragraphs.Item($Index)\n    $range = $para.Range\n    # Range.Text for a paragraph includes the trailing paragraph mark; strip it before comparing.\n    $current = $range.Text.TrimEnd([char]13)\n    if ($current -ne $Expected) {\n        throw \"Paragraph $Index text mismatch. Got: $current\"\n    }\n    $range.Text = $NewText\n}\n\n# Paragraph 3 (ShortSummary)\n$ShortSummaryOld = 'The article is about how ChatGPT and AI are expected to impact the Saudi workforce productivity by automating mundane tasks, providing tailor-made training programs, and reshaping various sectors like healthcare, transportation, energy, finance, and retail.'\n$ShortSummaryNew = 'The article is about the potential benefits of ChatGPT and AI on Saudi Arabia s workforce productivity, alleviating fears about job losses due to automation by highlighting their role in upskilling employees, boosting productivity, and creating new career opportunities.'\nSet-ParaText $d 3 $ShortSummaryOld $ShortSummaryNew\n\n# Paragraph 5 (Summary)\n$SummaryOld = '1. ChatGPT and AI can boost productivity in Saudi Arabia by providing tailor-made training programs, access to customized online courses, and fostering collaboration among team members.' + $br + '2. Embedding AI into operations requires a holistic approach that defines strategic objectives, advantages, and disadvantages. Understanding operational bottlenecks and how AI can address them is crucial.' + $br + '3. The automation of tasks by ChatGPT may lead to job displacement for employees entrusted with mundane and repetitive tasks, necessitating retraining or upskilling for workers.' + $br + '4. A positive employee and client experience is vital when implementing AI. Leadership, a futuristic view, and agility in making timely changes are essential for successful integration of AI into organizations.'\n$SummaryNew = '1. ChatGPT and AI can positively impact Saudi workforce productivity by providing tailor-made training programs, access to customized online courses, fostering collaboration and communication among team members.' + $br + '2. Embracing AI will require a strategic approach that defines objectives, advantages, and disadvantages, with an emphasis on addressing operational bottlenecks and choosing the optimal AI tools for the organization''s needs.' + $br + '3. The automation of tasks by ChatGPT might replace workers entrusted with mundane and repetitive functions, necessitating employees to acquire new skills through retraining or upskilling.' + $br + '4. Implementing AI in Saudi companies will require a holistic approach, fostering an optimistic environment for learning and improvement, solid leadership, and agility in making timely changes as needed.'\nSet-ParaText $d 5 $SummaryOld $SummaryNew\n\n# Paragraph 7 (Question1)\n$Question1Old = 'In the given text, the media frames the public discussion about ChatGPT and AI in general as a tool for boosting productivity and economic development. The metaphor of a wave of change is used to describe the impact of these technologies on the global workforce. Another recurring metaphor is that of an opportunity for organizations and employees to innovate, learn, collaborate, and improve their skills to adapt to an AI-enabled digital world. The text also refers to AI as a tool that can help automate mundane and repetitive tasks to free up time for strategic activities, using the metaphor of a productive and efficient workforce. However, there is also a concern about job displacement due to automation, which is depicted as a potential negative side effect. Overall, the text presents AI, particularly ChatGPT, as a transformative technology that can reshape and emancipate the Saudi workforce, but requires careful management and strategic planning to ensure positive outcomes.'\n$Question1New = 'In the text, the public discussion about ChatGPT and AI is framed as a tool for boosting productivity in Saudi Arabia''s workforce. The metaphors used to describe this process include:' + $br + '' + $br + '1. Wave of change (implying the transformative impact of AI technologies on the global workforce)' + $br + '2. Opportunity to boost productivity (indicating the potential positive outcomes of embracing innovation and AI)' + $br + '3. Embracing innovation (emphasizing the need for organizations to adapt and utilize AI platforms)' + $br + '4. Nurturing the right talent with a strong AI culture (implying the importance of fostering a supportive environment for AI implementation)' + $br + '5. Reshaping and even emancipating workforce (suggesting that AI can lead to significant changes in the workforce, potentially empowering it)' + $br + '6. Tailor-made training programs (indicating personalized learning opportunities provided by AI platforms like ChatGPT)' + $br + '7. Automation of tasks (implying the removal of mundane and repetitive tasks from employees)' + $br + '8. Acquiring new skills through retraining or upskilling (suggesting that workers will need to adapt to a more AI-enabled digital world)' + $br + '9. AI-enabled digital world (emphasizing the growing reliance on artificial intelligence in various sectors)' + $br + '10. Strategic roadmap (implying a clear plan for integrating AI into organizations'' operations)' + $br + '11. Agility (emphasizing the need for adaptability and quick changes as technological advancements occur)' + $br + '' + $br + 'These metaphors suggest that the media presents ChatGPT and AI as a transformative force, requiring a strategic approach to integration, and necessitating employee adaptation and skill acquisition in order to thrive in an increasingly AI-dependent world.'\nSet-ParaText $d 7 $Question1Old $Question1New\n\n# Paragraph 9 (Question2)\n$Question2Old = 'Perspectives and aspects that are being widely covered in this text include  1. The positive impact of AI technologies such as ChatGPT on Saudi Arabian workforce productivity. 2. Fear among workers about losing their jobs due to AI, and efforts to alleviate these fears by nurturing the right talent with a strong AI culture. 3. The role of AI in various sectors such as healthcare, transportation, energy, finance, and retail. 4. The benefits of AI for recruitment, hiring, training, development, upskilling, reskilling, talent collaboration, and knowledge management. 5. The potential of AI to increase productivity by automating mundane tasks and freeing up employees to focus on strategic activities. 6. The need to create a positive employee experience through understanding the strategic objectives, advantages, and disadvantages of AI implementation. 7. The importance of leadership, agility, and an optimistic environment for learning and improvement in embedding AI into organizations. 8. The potential impact of AI on Saudi Arabia s public sector operations and service delivery. 9. The concerns about the stability of worker employment due to automation. 10. The need for employees to acquire new skills through retraining or upskilling to remain marketable in an increasingly digital world. Aspects that are being ignored in this text include  1. The specific challenges faced by Saudi Arabia in implementing AI technologies, such as infrastructure limitations, regulatory constraints, or cultural resistance. 2. The ethical considerations related to the use of AI, including issues of data privacy, bias, and accountability. 3. The potential negative consequences of widespread AI adoption on employment, income distribution, and social inequality. 4. The environmental impact of AI technologies, such as energy consumption, e-waste generation, or resource depletion. 5. The potential risks associated with over-reliance on AI, such as algorithmic failures, cybersecurity threats, or dependence on a single vendor. 6. The need for continuous monitoring and evaluation of AI systems to ensure their performance remains optimal and aligned with strategic goals. 7. The need for international cooperation in the development and regulation of AI technologies to prevent a digital divide between developed and developing countries.'\n$Question2New = 'The text mainly covers the positive impact of ChatGPT and AI on Saudi workforce productivity, employee training, recruitment, and collaboration, as well as their potential applications in various sectors such as healthcare, transportation, energy, finance, and retail. It also discusses the need to nurture a strong AI culture within organizations, emphasizing innovation, experimentation, learning, and collaboration. The text mentions the concerns about job displacement due to AI but focuses more on how employees can acquire new skills through retraining or upskilling to adapt to an AI-enabled work environment.' + $br + '' + $br + 'Aspects that seem to be ignored in this text include the potential negative social impacts of widespread AI adoption, such as increased income inequality and loss of privacy, the ethical considerations around AI decision-making, and the need for regulations and guidelines to ensure fair use of AI technologies. Additionally, the text does not discuss the possible challenges or resistance from employees or organizations in adopting AI, nor does it delve into the long-term implications of AI on the Saudi economy beyond productivity improvements.'\nSet-ParaText $d 9 $Question2Old $Question2New\n\n# Paragraph 11 (Question3)\n$Question3Old = 'Not mentioned'\n$Question3New = 'The Arabic world is not explicitly mentioned in the text, but it can be inferred that the Kingdom (presumably Saudi Arabia) is leveraging AI to boost economic development and increase workforce productivity. The AI platform ChatGPT is highlighted as a useful tool for this purpose. The article discusses how AI, specifically ChatGPT, can impact recruitment, hiring, training, upskilling, reskilling, talent collaboration, knowledge management, and various sectors like healthcare, transportation, energy, finance, and retail. It also mentions the importance of nurturing a strong AI culture to maximize the benefits of AI technologies.'\nSet-ParaText $d 11 $Question3Old $Question3New\n\n# Paragraph 13 (Question4)\n$Question4Old = 'The article suggests that AI, specifically ChatGPT, can positively impact Saudi workforce productivity by providing customized training programs, collaboration tools, and automating mundane tasks. However, it also warns about potential job displacement due to automation. Embedding AI requires a holistic approach, clear strategic objectives, understanding of operational bottlenecks, and fostering an optimistic learning environment. Leadership, agility, and a futuristic view are crucial for successful implementation.'\n$Question4New = 'The article suggests that ChatGPT and AI can boost Saudi Arabia''s workforce productivity by automating mundane tasks and providing tailored training programs. However, it emphasizes the need to foster an optimistic learning environment and nurture a strong AI culture to address employees'' concerns about job displacement. Embracing innovation can lead to increased efficiency, better positions for employees, and growth opportunities in various sectors like healthcare, transportation, finance, and retail.'\nSet-ParaText $d 13 $Question4Old $Question4New\n\n# Paragraph 17 (Entities)\n$EntitiesOld = 'Embedding AI, Arthur D. Little, Saudi Arabias, Kingdoms, Kaspersky, Arab News, Raymond Khoury, Khoury, AI, Khourys, Kasperskys, Arthur D. Little'\n$EntitiesNew = 'Kaspersky, AI, Raymond Khoury, Kasperskys, Embedding AI, Kingdoms, Arthur D. Little, Arab News, Khourys, Arthur D. Little, Saudi Arabias, Khoury'\nSet-ParaText $d 17 $EntitiesOld $EntitiesNew\n"}
